# replace buffer.name with buffer.id
# The "buffer" sheet (sheet1) had a leading "name" column that duplicated
# information already present via the other columns; drop it and let the
# remaining columns (type, location, item, onhand) shift left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buffer")

# Remove the whole "name" column (column A). This shifts type/location/item/onhand
# one column to the left, exactly like Excel's Delete Entire Column command.
$ws.Columns.Item(1).Delete()

# Make "buffer" the active sheet and select column A (now "type"), matching
# the end-user state captured after the edit.
$ws.Activate()
$ws.Columns.Item(1).Select()
